$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Mark every currently-"Pendente" requirement in the Status column (A) as
# "Feito" (done). Rows 6-9 and 11-20 are "Pendente"; rows 3-5 and 10 are
# already "Feito" and are left untouched. Once no row is left with
# "Pendente", the shared string is dropped automatically when the sheet is
# re-serialized.
$pendingRows = @(6, 7, 8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
foreach ($r in $pendingRows) {
    $ws.Range("A$r").Value = "Feito"
}

# Force recalculation so the COUNTIF summary cells (K3 = "feito" count,
# K4 = "pendente" count) pick up the new totals.
$excel.CalculateFullRebuild()

# Activate the "Product Backlog" sheet (it becomes the saved active tab
# instead of "5W2H") and move the selection/scroll position to A2.
$ws.Activate()
$ws.Range("A2").Select()
